$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.847.22'
$ws.Range("E2").Value = '  -1.90%  '

$ws.Range("D3").Value = '3.512.26'
$ws.Range("E3").Value = '  -3.59%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.86%  '

$ws.Range("E7").Value = '  -2.55%  '

$ws.Range("D8").Value = '3.501.03'
$ws.Range("E8").Value = '  -3.51%  '

$ws.Range("E10").Value = '  -6.53%  '

$ws.Range("E11").Value = '  -4.33%  '

$ws.Range("E12").Value = '  -4.90%  '

$ws.Range("E13").Value = '  -6.61%  '

$ws.Range("E14").Value = '  -4.61%  '

$ws.Range("D15").Value = '4.064.75'
$ws.Range("E15").Value = '  -3.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '651.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.24%  '

$ws.Range("D17").Value = '69.859.20'
$ws.Range("E17").Value = '  -2.01%  '

$ws.Range("D18").Value = '3.513.50'
$ws.Range("E18").Value = '  -3.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.64%  '

$ws.Range("E21").Value = '  -1.95%  '

$ws.Range("E22").Value = '  -5.00%  '

$ws.Range("E23").Value = '  -3.75%  '

$ws.Range("E24").Value = '  -0.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.14%  '

$ws.Range("E27").Value = '  -4.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.78'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.53%  '

$ws.Range("E31").Value = '  -6.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.86%  '

$ws.Range("E33").Value = '  -5.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '579.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.06%  '

$ws.Range("E35").Value = '  -4.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '61.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.21%  '

$ws.Range("D37").Value = '3.760.30'
$ws.Range("E37").Value = '  -4.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("E39").Value = '  -9.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +44.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.89'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.02%  '

$ws.Range("E43").Value = '  -4.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.134'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.84%  '

$ws.Range("E46").Value = '  -3.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.59%  '

$ws.Range("E48").Value = '  -3.83%  '

$ws.Range("E49").Value = '  -3.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.38%  '
